# Applies updated market-price derived values to the leve profit sheets
# (mirrors a scheduled data refresh across all 8 job sheets)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 753.73334
$ws.Range("I98").Value = 557.9167
$ws.Range("J98").Value = 1537
$ws.Range("K98").Value = 557.9167
$ws.Range("L98").Value = 1537
$ws.Range("M98").Value = 940.0833
$ws.Range("N98").Value = -4533
$ws.Range("H116").Value = 5299.615
$ws.Range("I116").Value = 2899.8
$ws.Range("K116").Value = 2899.8
$ws.Range("M116").Value = 542.1999999999998
$ws.Range("H122").Value = 753.73334
$ws.Range("I122").Value = 557.9167
$ws.Range("J122").Value = 1537
$ws.Range("K122").Value = 1673.7501
$ws.Range("L122").Value = 4611
$ws.Range("M122").Value = 776.2499
$ws.Range("N122").Value = -9511
$ws.Range("H138").Value = 2496.1667
$ws.Range("I138").Value = 1529.4375
$ws.Range("J138").Value = 3269.55
$ws.Range("K138").Value = 4588.3125
$ws.Range("L138").Value = 9808.650000000001
$ws.Range("M138").Value = 551.6875
$ws.Range("N138").Value = -20088.65

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 657.4839
$ws.Range("I2").Value = 690.4167
$ws.Range("K2").Value = 690.4167
$ws.Range("M2").Value = -577.4167
$ws.Range("H32").Value = 1977.8
$ws.Range("I32").Value = 1528.7333
$ws.Range("J32").Value = 6019.4
$ws.Range("K32").Value = 1528.7333
$ws.Range("L32").Value = 6019.4
$ws.Range("M32").Value = -1241.7333
$ws.Range("N32").Value = -6593.4
$ws.Range("H45").Value = 2452.4443
$ws.Range("I45").Value = 1828.2858
$ws.Range("J45").Value = 3124.6155
$ws.Range("K45").Value = 1828.2858
$ws.Range("L45").Value = 3124.6155
$ws.Range("M45").Value = -1451.2858
$ws.Range("N45").Value = -3878.6155
$ws.Range("H74").Value = 41668656
$ws.Range("I74").Value = 76923730
$ws.Range("J74").Value = 3570.818
$ws.Range("K74").Value = 76923730
$ws.Range("L74").Value = 3570.818
$ws.Range("M74").Value = -76922856
$ws.Range("N74").Value = -5318.818
$ws.Range("H77").Value = 41668656
$ws.Range("I77").Value = 76923730
$ws.Range("J77").Value = 3570.818
$ws.Range("K77").Value = 384618650
$ws.Range("L77").Value = 17854.09
$ws.Range("M77").Value = -384614282
$ws.Range("N77").Value = -26590.09
$ws.Range("H116").Value = 657.4839
$ws.Range("I116").Value = 690.4167
$ws.Range("K116").Value = 690.4167
$ws.Range("M116").Value = 1603.5833
$ws.Range("H132").Value = 13269.954
$ws.Range("I132").Value = 1825.3143
$ws.Range("J132").Value = 57776.89
$ws.Range("K132").Value = 5475.9429
$ws.Range("L132").Value = 173330.67
$ws.Range("M132").Value = -2945.9429
$ws.Range("N132").Value = -178390.67

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 657.4839
$ws.Range("I3").Value = 690.4167
$ws.Range("K3").Value = 690.4167
$ws.Range("M3").Value = -576.4167
$ws.Range("H22").Value = 710.4286
$ws.Range("I22").Value = 660.6667
$ws.Range("K22").Value = 660.6667
$ws.Range("M22").Value = -487.6667
$ws.Range("H105").Value = 1771.6327
$ws.Range("I105").Value = 1436.4286
$ws.Range("K105").Value = 1436.4286
$ws.Range("M105").Value = 310.5714
$ws.Range("H107").Value = 975.4
$ws.Range("I107").Value = 835.9167
$ws.Range("K107").Value = 835.9167
$ws.Range("M107").Value = 1084.0833

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 50004070
$ws.Range("J62").Value = 5234.1665
$ws.Range("L62").Value = 5234.1665
$ws.Range("N62").Value = -6482.1665
$ws.Range("H65").Value = 50004070
$ws.Range("J65").Value = 5234.1665
$ws.Range("L65").Value = 26170.8325
$ws.Range("N65").Value = -32410.8325
$ws.Range("H99").Value = 23812958
$ws.Range("I99").Value = 2935.7144
$ws.Range("J99").Value = 71433000
$ws.Range("K99").Value = 2935.7144
$ws.Range("L99").Value = 71433000
$ws.Range("M99").Value = -1437.7144
$ws.Range("N99").Value = -71435996
$ws.Range("H126").Value = 23812958
$ws.Range("I126").Value = 2935.7144
$ws.Range("J126").Value = 71433000
$ws.Range("K126").Value = 8807.143199999999
$ws.Range("L126").Value = 214299000
$ws.Range("M126").Value = -6337.143199999999
$ws.Range("N126").Value = -214303940
$ws.Range("H132").Value = 3776.3157
$ws.Range("I132").Value = 2754.6155
$ws.Range("J132").Value = 5990
$ws.Range("K132").Value = 8263.8465
$ws.Range("L132").Value = 17970
$ws.Range("M132").Value = -5733.8465
$ws.Range("N132").Value = -23030
$ws.Range("H134").Value = 1265.3334
$ws.Range("I134").Value = 1109.9286
$ws.Range("K134").Value = 3329.7858
$ws.Range("M134").Value = -794.7857999999997

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H106").Value = 4088.1667
$ws.Range("J106").Value = 4088.1667
$ws.Range("L106").Value = 12264.5001
$ws.Range("N106").Value = -14156.5001
$ws.Range("H121").Value = 2446.2
$ws.Range("J121").Value = 2982.75
$ws.Range("L121").Value = 8948.25
$ws.Range("N121").Value = -11568.25
$ws.Range("H127").Value = 695.55554
$ws.Range("J127").Value = 695.55554
$ws.Range("L127").Value = 2086.66662
$ws.Range("N127").Value = -12006.66662

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 13236.182
$ws.Range("I70").Value = 5497.5
$ws.Range("J70").Value = 17658.285
$ws.Range("K70").Value = 5497.5
$ws.Range("L70").Value = 17658.285
$ws.Range("M70").Value = -5227.5
$ws.Range("N70").Value = -18198.285
$ws.Range("H73").Value = 13236.182
$ws.Range("I73").Value = 5497.5
$ws.Range("J73").Value = 17658.285
$ws.Range("K73").Value = 5497.5
$ws.Range("L73").Value = 17658.285
$ws.Range("M73").Value = -4561.5
$ws.Range("N73").Value = -19530.285
$ws.Range("H80").Value = 3816.2856
$ws.Range("I80").Value = 3313.5715
$ws.Range("J80").Value = 4067.6428
$ws.Range("K80").Value = 3313.5715
$ws.Range("L80").Value = 4067.6428
$ws.Range("M80").Value = -2315.5715
$ws.Range("N80").Value = -6063.6428
$ws.Range("H83").Value = 3816.2856
$ws.Range("I83").Value = 3313.5715
$ws.Range("J83").Value = 4067.6428
$ws.Range("K83").Value = 16567.8575
$ws.Range("L83").Value = 20338.214
$ws.Range("M83").Value = -11575.8575
$ws.Range("N83").Value = -30322.214
$ws.Range("H102").Value = 20001644
$ws.Range("I102").Value = 22728582
$ws.Range("J102").Value = 4102.6665
$ws.Range("K102").Value = 22728582
$ws.Range("L102").Value = 4102.6665
$ws.Range("M102").Value = -22726960
$ws.Range("N102").Value = -7346.6665
$ws.Range("H113").Value = 3612.919
$ws.Range("J113").Value = 1642.1666
$ws.Range("L113").Value = 1642.1666
$ws.Range("N113").Value = -5982.1666
$ws.Range("H126").Value = 3091.8667
$ws.Range("I126").Value = 2174.3809
$ws.Range("J126").Value = 5232.6665
$ws.Range("K126").Value = 6523.1427
$ws.Range("L126").Value = 15697.9995
$ws.Range("M126").Value = -4053.1427
$ws.Range("N126").Value = -20637.9995
$ws.Range("H132").Value = 24222.8
$ws.Range("I132").Value = 4383.6
$ws.Range("K132").Value = 13150.8
$ws.Range("M132").Value = -10620.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4407
$ws.Range("I7").Value = 3962.375
$ws.Range("J7").Value = 4999.8335
$ws.Range("K7").Value = 3962.375
$ws.Range("L7").Value = 4999.8335
$ws.Range("M7").Value = -3850.375
$ws.Range("N7").Value = -5223.8335
$ws.Range("H40").Value = 4991.4
$ws.Range("I40").Value = 3326.6667
$ws.Range("J40").Value = 7488.5
$ws.Range("K40").Value = 3326.6667
$ws.Range("L40").Value = 7488.5
$ws.Range("M40").Value = -3190.6667
$ws.Range("N40").Value = -7760.5
$ws.Range("H126").Value = 4407
$ws.Range("I126").Value = 3962.375
$ws.Range("J126").Value = 4999.8335
$ws.Range("K126").Value = 11887.125
$ws.Range("L126").Value = 14999.5005
$ws.Range("M126").Value = -9417.125
$ws.Range("N126").Value = -19939.5005
$ws.Range("H132").Value = 367136.3
$ws.Range("I132").Value = 525205.2
$ws.Range("J132").Value = 3577.8
$ws.Range("K132").Value = 1575615.6
$ws.Range("L132").Value = 10733.4
$ws.Range("M132").Value = -1573085.6
$ws.Range("N132").Value = -15793.4
$ws.Range("H136").Value = 1645.909
$ws.Range("I136").Value = 1567.2222
$ws.Range("J136").Value = 2000
$ws.Range("K136").Value = 4701.6666
$ws.Range("L136").Value = 6000
$ws.Range("M136").Value = -2151.6666
$ws.Range("N136").Value = -11100

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H119").Value = 35000
$ws.Range("J119").Value = 35000
$ws.Range("L119").Value = 35000
$ws.Range("N119").Value = -44676
$ws.Range("H126").Value = 1455
$ws.Range("I126").Value = 1372.4
$ws.Range("J126").Value = 1609.875
$ws.Range("K126").Value = 4117.200000000001
$ws.Range("L126").Value = 4829.625
$ws.Range("M126").Value = -1647.200000000001
$ws.Range("N126").Value = -9769.625

Write-Host "Applied 225 cell updates across 8 sheets"
